# Update the two battery component names so that the dynamic parsing logic
# (which presumably splits on commas) is no longer tripped up by the commas
# inside the component descriptions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A30").Value = "CGAS007 Energizer Single Cell 3.7 V Lithium Ion 1 Ah"
$ws.Range("A31").Value = "CA5L Energizer Single Cell 3.7 V Lithium Ion 980 mAh"

# Reflect the author's final view/selection state: scrolled so row 22 is at
# the top and the active cell/selection is A30.
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A30").Select()
